$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.384.86"
$ws.Range("D3").Value = "2.512.22"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "539.46"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "139.52"
$ws.Range("E6").Value = "  -3.95%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "0.562"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").Value = "2.514.05"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "5.38"
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").Value = "0.358"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").Value = "2.961.53"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "23.40"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "59.280.40"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "0.0000141"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "2.512.79"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "11.10"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").Value = "4.30"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "324.77"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "63.22"
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("D25").Value = "0.423"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "7.82"
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("D29").Value = "6.86"
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("D30").Value = "0.0₃0779"
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").Value = "1.79"
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").Value = "164.36"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "1.45"
$ws.Range("E34").Value = "  -2.10%  "
$ws.Range("E35").Value = "  -7.17%  "
$ws.Range("D36").Value = "18.50"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("D37").Value = "4.26"
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("D38").Value = "1.59"
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "3.68"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  -3.03%  "
$ws.Range("D42").Value = "5.23"
$ws.Range("E42").Value = "  -6.31%  "
$ws.Range("D43").Value = "281.42"
$ws.Range("E43").Value = "  -4.72%  "
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").Value = "0.598"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").Value = "0.0936"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").Value = "124.27"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "0.0225"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "17.87"
$ws.Range("E51").Value = "  -2.10%  "
